# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 24; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 26; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 28; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 30; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 38; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 39; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 43; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 49; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 54; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 55; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 56; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 61; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 83; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 100; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 119; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 134; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 135; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 152; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 155; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 156; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 194; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 207; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 214; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 221; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 223; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 235; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 237; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 241; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 256; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 274; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 279; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 340; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 353; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 355; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 357; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 364; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 366; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 378; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 384; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 385; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 388; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 397; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 402; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 420; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 447; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 448; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 457; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 458; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 468; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 482; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 518; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 523; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 525; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 527; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 532; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 535; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 539; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 541; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 544; DAMSLTag = "ba"; DialogAct = "Appreciation" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows."
